# Applies the "Actualización automática" update: a new sale of 2925.73
# (PORCELANATO) is recorded for client TRUJILLO TORRES VINICIO RUBEN
# (advisor HIDALGO HIDALGO PEDRO GUSTAVO), and all of the dependent
# totals / ratios across the three sheets are refreshed accordingly.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" -------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M20").Value = 2925.73
$ws1.Range("M22").Value = "8 de 20"

# --- Sheet "VENTA MENSUAL" ----------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F20").Value = 2925.73
$ws2.Range("F22").Value = 40549.43

# --- Sheet "CUMPLIMIENTO MENSUAL" ---------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 16 - PORCELANATO group totals
$ws3.Range("D16").Value = 31336.38
$ws3.Range("E16").Value = 12929.86
$ws3.Range("F16").Value = 0.7079069738021572

# Row 19 - TOTAL row
$ws3.Range("D19").Value = 40549.43
$ws3.Range("E19").Value = 24828.56762291769
$ws3.Range("F19").Value = 0.6202305282256879
